$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: "Hoja1" -> "PatientStats"
$ws.Name = "PatientStats"

# --- Structural column changes (performed rightmost-first so earlier column
#     letters remain valid reference points for subsequent operations) ---

# Old column M (" WR Total Time") is dropped/consolidated away.
$ws.Columns.Item(13).Delete()

# A new blank column is inserted where old column F (" EX Timestamp") was,
# making room for the new " Call Back Timestamp" field.
$ws.Columns.Item(6).Insert()

# A new blank column is inserted where old column D (" Appt Time") was,
# making room for the new " Appt Type" field.
$ws.Columns.Item(4).Insert()

# --- Row 1 header relabeling / new fields ---
$ws.Range("A1").Value = "MRN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = " Physician"
$ws.Range("D1").Value = " Appt Type"
$ws.Range("E1").Value = " Appt Time"
$ws.Range("F1").Value = " Register Timestamp"
$ws.Range("G1").Value = " Call Back Timestamp"
$ws.Range("H1").Value = "Discharge Timestamp"
$ws.Range("I1").Value = " FC Start"
$ws.Range("J1").Value = " FC End"
$ws.Range("K1").Value = "Imaging Start Timestamp"
$ws.Range("L1").Value = "Imaging End Timestamp"
$ws.Range("M1").Value = "Imaging Total Time"
$ws.Range("N1").Value = "Callback vs. Appt Time"
$ws.Range("O1").Value = " EX Total Time"
$ws.Range("P1").Value = "Total Time"
$ws.Range("Q1").Value = " AT Entry"

# --- Selection / view ---
$ws.Range("O7").Select()

# --- Column widths (best effort; runtime rounds ColumnWidth to the nearest 1/6
#     character and adds an internal 5/6 padding offset, so we pre-compensate by
#     subtracting 5/6 from every desired target width) ---
$ws.Columns.Item(3).ColumnWidth  = 13.330729166666666   # C
$ws.Columns.Item(4).ColumnWidth  = 13.330729166666666   # D
$ws.Columns.Item(6).ColumnWidth  = 15.666666666666666   # F
$ws.Columns.Item(7).ColumnWidth  = 16.166666666666668   # G
$ws.Columns.Item(8).ColumnWidth  = 16.166666666666668   # H
$ws.Columns.Item(9).ColumnWidth  = 10.998697916666666   # I
$ws.Columns.Item(10).ColumnWidth = 11.166666666666666   # J
$ws.Columns.Item(11).ColumnWidth = 19.166666666666668   # K
$ws.Columns.Item(12).ColumnWidth = 19.166666666666668   # L
$ws.Columns.Item(13).ColumnWidth = 19.166666666666668   # M
$ws.Columns.Item(14).ColumnWidth = 17.166666666666668   # N
$ws.Columns.Item(15).ColumnWidth = 13.330729166666666   # O
$ws.Columns.Item(16).ColumnWidth = 12.498697916666666   # P
$ws.Columns.Item(17).ColumnWidth = 11.830729166666666   # Q

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
